$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Insert a new row at position 6 (shifts old rows 6,7,8 down to 7,8,9;
#    also auto-extends the shared formula in column E and the SUM() range).
$ws.Rows.Item(6).Insert()

# 2. Populate the new row (Digital Temp / Bluetooth-adjacent sensor part).
#    Write the URL first, then the part name, so the new shared-string
#    entries land in the same order as the target workbook.
$ws.Range("F6").Value2 = "https://www.sparkfun.com/products/11295"
$ws.Range("A6").Value2 = "Digital Temp"

# Give F6 the same "looks like a link but isn't wired up" styling the GPS
# row (F5) uses - centered/underlined Hyperlink look without an actual
# Hyperlinks collection entry.
$ws.Range("F6").Style = "Hyperlink"

# 3. The row insert does not move the worksheet's Hyperlinks collection
#    along with the shifted cells, so rebuild it from scratch in the
#    correct final order/positions. (Deleting any one hyperlink clears the
#    whole collection in this host, so we delete once and re-add all five.)
$ws.Range("F2").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("F2"), "http://www.newark.com/stmicroelectronics/ld1117s33ctr/ic-ldo-volt-reg-3-3v-0-8a-sot/dp/89K0626?CMP=AFC-OP") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.sparkfun.com/products/9473") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "http://www.newark.com/nxp/74hc4052d-653/ic-analog-mux-dmux-dual-4-x-1/dp/78R7402") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.sparkfun.com/products/9609") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), "http://www.mouser.com/ProductDetail/Texas-Instruments/EK-LM4F120XL/?qs=t9Lg9qrXjEy2enepSwqR9A==") | Out-Null

# Re-adding a hyperlink mints a fresh cell style; snap those five cells
# back onto the workbook's normal Hyperlink style so no stray formatting
# differences are introduced.
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("F7").Style = "Hyperlink"
$ws.Range("F8").Style = "Hyperlink"

# 4. Match the saved selection state from the edited workbook.
$ws.Range("F11").Select() | Out-Null
